# Auto-generated script: refresh market-price derived columns (H-N)
# on the Leve profit tracker sheets, per scheduled-runner data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 566.3333
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 599.5
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 599.5
$ws.Range("M19").Value = -325
$ws.Range("N19").Value = -949.5

$ws.Range("H33").Value = 300.33334
$ws.Range("I33").Value = 307.5
$ws.Range("K33").Value = 307.5
$ws.Range("M33").Value = -78.5

$ws.Range("H116").Value = 3698.2
$ws.Range("I116").Value = 2630
$ws.Range("J116").Value = 4410.3335
$ws.Range("K116").Value = 2630
$ws.Range("L116").Value = 4410.3335
$ws.Range("M116").Value = 812
$ws.Range("N116").Value = -11294.3335

$ws.Range("H132").Value = 3655.6538
$ws.Range("I132").Value = 4220.7144
$ws.Range("K132").Value = 12662.1432
$ws.Range("M132").Value = -10132.1432

$ws.Range("H137").Value = 1931.25
$ws.Range("I137").Value = 2040
$ws.Range("K137").Value = 6120
$ws.Range("M137").Value = -3570


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 657.9
$ws.Range("I2").Value = 716.72
$ws.Range("J2").Value = 363.8
$ws.Range("K2").Value = 716.72
$ws.Range("L2").Value = 363.8
$ws.Range("M2").Value = -603.72
$ws.Range("N2").Value = -589.8

$ws.Range("H32").Value = 4920.5303
$ws.Range("I32").Value = 4177.493
$ws.Range("J32").Value = 10344.7
$ws.Range("K32").Value = 4177.493
$ws.Range("L32").Value = 10344.7
$ws.Range("M32").Value = -3890.493
$ws.Range("N32").Value = -10918.7

$ws.Range("H61").Value = 1745.4546
$ws.Range("I61").Value = 1665
$ws.Range("J61").Value = 2550
$ws.Range("K61").Value = 1665
$ws.Range("L61").Value = 2550
$ws.Range("M61").Value = -1453
$ws.Range("N61").Value = -2974

$ws.Range("H63").Value = 2031.8572
$ws.Range("J63").Value = 1899
$ws.Range("L63").Value = 1899
$ws.Range("N63").Value = -3271

$ws.Range("H66").Value = 2031.8572
$ws.Range("J66").Value = 1899
$ws.Range("L66").Value = 9495
$ws.Range("N66").Value = -16359

$ws.Range("H74").Value = 62502260
$ws.Range("I74").Value = 200000640
$ws.Range("J74").Value = 2995.4546
$ws.Range("K74").Value = 200000640
$ws.Range("L74").Value = 2995.4546
$ws.Range("M74").Value = -199999766
$ws.Range("N74").Value = -4743.4546

$ws.Range("H77").Value = 62502260
$ws.Range("I77").Value = 200000640
$ws.Range("J77").Value = 2995.4546
$ws.Range("K77").Value = 1000003200
$ws.Range("L77").Value = 14977.273
$ws.Range("M77").Value = -999998832
$ws.Range("N77").Value = -23713.273

$ws.Range("H102").Value = 1179.091
$ws.Range("I102").Value = 1133.75
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 1133.75
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = 488.25
$ws.Range("N102").Value = -4544

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

$ws.Range("H110").Value = 721.2857
$ws.Range("I110").Value = 633.3333
$ws.Range("J110").Value = 879.6
$ws.Range("K110").Value = 633.3333
$ws.Range("L110").Value = 879.6
$ws.Range("M110").Value = 1411.6667
$ws.Range("N110").Value = -4969.6

$ws.Range("H112").Value = 34837.6
$ws.Range("J112").Value = 34837.6
$ws.Range("L112").Value = 34837.6
$ws.Range("N112").Value = -37791.6

$ws.Range("H116").Value = 657.9
$ws.Range("I116").Value = 716.72
$ws.Range("J116").Value = 363.8
$ws.Range("K116").Value = 716.72
$ws.Range("L116").Value = 363.8
$ws.Range("M116").Value = 1577.28
$ws.Range("N116").Value = -4951.8

$ws.Range("H122").Value = 1752.9688
$ws.Range("I122").Value = 1559.5
$ws.Range("K122").Value = 4678.5
$ws.Range("M122").Value = -2228.5

$ws.Range("H136").Value = 1745.4546
$ws.Range("I136").Value = 1665
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 4995
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -2445
$ws.Range("N136").Value = -12750


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 657.9
$ws.Range("I3").Value = 716.72
$ws.Range("J3").Value = 363.8
$ws.Range("K3").Value = 716.72
$ws.Range("L3").Value = 363.8
$ws.Range("M3").Value = -602.72
$ws.Range("N3").Value = -591.8

$ws.Range("H94").Value = 940.4666999999999
$ws.Range("I94").Value = 830.7
$ws.Range("K94").Value = 830.7
$ws.Range("M94").Value = -379.7

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = ""

$ws.Range("H107").Value = 717.5454999999999
$ws.Range("I107").Value = 738
$ws.Range("K107").Value = 738
$ws.Range("M107").Value = 1182


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3763.1943
$ws.Range("I31").Value = 3966
$ws.Range("J31").Value = 3695.5925
$ws.Range("K31").Value = 3966
$ws.Range("L31").Value = 3695.5925
$ws.Range("M31").Value = -3671
$ws.Range("N31").Value = -4285.592500000001

$ws.Range("H34").Value = 3763.1943
$ws.Range("I34").Value = 3966
$ws.Range("J34").Value = 3695.5925
$ws.Range("K34").Value = 3966
$ws.Range("L34").Value = 3695.5925
$ws.Range("M34").Value = -3764
$ws.Range("N34").Value = -4099.592500000001

$ws.Range("H92").Value = 21999
$ws.Range("J92").Value = 21999
$ws.Range("L92").Value = 21999
$ws.Range("N92").Value = -26991

$ws.Range("H99").Value = 3844.1177
$ws.Range("I99").Value = 2950
$ws.Range("J99").Value = 6750
$ws.Range("K99").Value = 2950
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -1452
$ws.Range("N99").Value = -9746

$ws.Range("H126").Value = 3844.1177
$ws.Range("I126").Value = 2950
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 8850
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -6380
$ws.Range("N126").Value = -25190

$ws.Range("H141").Value = 121109.75
$ws.Range("J141").Value = 121109.75
$ws.Range("L141").Value = 121109.75
$ws.Range("N141").Value = -131469.75


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 754.375
$ws.Range("J122").Value = 836.3077
$ws.Range("L122").Value = 7526.7693
$ws.Range("N122").Value = -12426.7693

$ws.Range("H131").Value = 785.88043
$ws.Range("J131").Value = 780.45557
$ws.Range("L131").Value = 2341.36671
$ws.Range("N131").Value = -12421.36671


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3046
$ws.Range("I80").Value = 2603.2222
$ws.Range("J80").Value = 3235.762
$ws.Range("K80").Value = 2603.2222
$ws.Range("L80").Value = 3235.762
$ws.Range("M80").Value = -1605.2222
$ws.Range("N80").Value = -5231.762000000001

$ws.Range("H83").Value = 3046
$ws.Range("I83").Value = 2603.2222
$ws.Range("J83").Value = 3235.762
$ws.Range("K83").Value = 13016.111
$ws.Range("L83").Value = 16178.81
$ws.Range("M83").Value = -8024.111000000001
$ws.Range("N83").Value = -26162.81

$ws.Range("H123").Value = 6024.1816
$ws.Range("I123").Value = 3437.7778
$ws.Range("K123").Value = 3437.7778
$ws.Range("M123").Value = -987.7777999999998

$ws.Range("H126").Value = 4633.3335
$ws.Range("I126").Value = 2850
$ws.Range("J126").Value = 6671.4287
$ws.Range("K126").Value = 8550
$ws.Range("L126").Value = 20014.2861
$ws.Range("M126").Value = -6080
$ws.Range("N126").Value = -24954.2861


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1281.6471
$ws.Range("I46").Value = 1034.8572
$ws.Range("J46").Value = 2433.3333
$ws.Range("K46").Value = 1034.8572
$ws.Range("L46").Value = 2433.3333
$ws.Range("M46").Value = -846.8571999999999
$ws.Range("N46").Value = -2809.3333

$ws.Range("H104").Value = 18642.777
$ws.Range("J104").Value = 18642.777
$ws.Range("L104").Value = 18642.777
$ws.Range("N104").Value = -25630.777


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14000
$ws.Range("J54").Value = 14000
$ws.Range("L54").Value = 14000
$ws.Range("N54").Value = -15040

$ws.Range("H101").Value = 13250.5
$ws.Range("J101").Value = 13250.5
$ws.Range("L101").Value = 13250.5
$ws.Range("N101").Value = -19740.5

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

$ws.Range("H113").Value = 1587.9231
$ws.Range("I113").Value = 2561.5715
$ws.Range("J113").Value = 452
$ws.Range("K113").Value = 7684.7145
$ws.Range("L113").Value = 1356
$ws.Range("M113").Value = -5514.7145
$ws.Range("N113").Value = -5696

$ws.Range("H126").Value = 1203.5217
$ws.Range("I126").Value = 1230.9546
$ws.Range("J126").Value = 600
$ws.Range("K126").Value = 3692.8638
$ws.Range("L126").Value = 1800
$ws.Range("M126").Value = -1222.8638
$ws.Range("N126").Value = -6740

$ws.Range("H132").Value = 1886.2632
$ws.Range("I132").Value = 1084.1
$ws.Range("J132").Value = 2777.5557
$ws.Range("K132").Value = 3252.3
$ws.Range("L132").Value = 8332.667099999999
$ws.Range("M132").Value = -722.2999999999997
$ws.Range("N132").Value = -13392.6671

